# Refresh the crypto price/volume table (cryptos.xlsx) with the latest
# scrape, mirroring the GitHub Actions update job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") holds plain numeric-looking text (e.g. "210.74").
# Excel auto-converts such literals to real Numbers on assignment, which
# would change the cell type away from the source Text cells. Prefixing
# with a leading apostrophe is Excel's own mechanism to force the literal
# to be stored as text (the apostrophe itself is not part of the value),
# so we use it only where the new price value would otherwise parse as a
# number.
function Set-CellText($range, $value) {
    if ($value -match "^\s*[+-]?\d+(\.\d+)?\s*$") {
        $ws.Range($range).Value = "'" + $value
    } else {
        $ws.Range($range).Value = $value
    }
}

Set-CellText "D2" '26.650.73'
Set-CellText "E2" '  +0.17%  '
Set-CellText "D3" '1.598.64'
Set-CellText "E3" '  +0.71%  '
Set-CellText "E4" '  -0.11%  '
Set-CellText "D5" '210.74'
Set-CellText "E5" '  -0.18%  '
Set-CellText "E6" '  +1.46%  '
Set-CellText "E7" '  -0.09%  '
Set-CellText "E8" '  -0.07%  '
Set-CellText "E9" '  -1.34%  '
Set-CellText "E10" '  +0.68%  '
Set-CellText "E11" '  +0.46%  '
Set-CellText "D12" '1.822.25'
Set-CellText "E12" '  +0.54%  '
Set-CellText "D13" '1.617.11'
Set-CellText "E13" '  +2.01%  '
Set-CellText "E14" '  -0.53%  '
Set-CellText "E15" '  -1.28%  '
Set-CellText "D16" '64.83'
Set-CellText "E16" '  +1.69%  '
Set-CellText "D17" '26.632.23'
Set-CellText "E17" '  -0.02%  '
Set-CellText "D18" '0.0₃0729'
Set-CellText "E18" '  -0.01%  '
Set-CellText "D19" '209.05'
Set-CellText "E19" '  +0.20%  '
Set-CellText "E20" '  +0.01%  '
Set-CellText "D21" '6.75'
Set-CellText "E21" '  +1.19%  '
Set-CellText "E22" '  +0.15%  '
Set-CellText "E23" '  -3.24%  '
Set-CellText "E24" '  +0.26%  '
Set-CellText "D25" '145.78'
Set-CellText "E25" '  -0.40%  '
Set-CellText "E26" '  -0.13%  '
Set-CellText "E27" '  -3.04%  '
Set-CellText "E28" '  +2.35%  '
Set-CellText "D29" '15.28'
Set-CellText "E29" '  -0.15%  '
Set-CellText "D30" '0.0505'
Set-CellText "E30" '  +0.99%  '
Set-CellText "E31" '  -0.22%  '
Set-CellText "E32" '  -0.37%  '
Set-CellText "D33" '0.661'
Set-CellText "E33" '  +1.50%  '
Set-CellText "E34" '  -0.25%  '
Set-CellText "D35" '1.294.02'
Set-CellText "E35" '  -0.78%  '
Set-CellText "E36" '  +0.34%  '
Set-CellText "E37" '  -1.28%  '
Set-CellText "E38" '  -0.60%  '
Set-CellText "E39" '  +2.97%  '
Set-CellText "E40" '  -0.05%  '
Set-CellText "B41" 'FraxShare'
Set-CellText "C41" 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-CellText "D41" '5.40'
Set-CellText "E41" '  +2.12%  '
Set-CellText "B42" 'MXToken'
Set-CellText "C42" 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-CellText "D42" '2.21'
Set-CellText "E42" '  +1.66%  '
Set-CellText "D43" '0.785'
Set-CellText "E43" '  -0.07%  '
Set-CellText "D44" '63.78'
Set-CellText "E44" '  +1.35%  '
Set-CellText "D45" '1.735.07'
Set-CellText "E45" '  +0.48%  '
Set-CellText "D46" '0.891'
Set-CellText "E46" '  +7.28%  '
Set-CellText "D47" '90.17'
Set-CellText "E47" '  +1.39%  '
Set-CellText "E48" '  +0.10%  '
Set-CellText "E49" '  +2.53%  '
Set-CellText "E50" '  -0.51%  '
Set-CellText "D51" '7.46'
Set-CellText "E51" '  -0.26%  '
